# Split two merged runs in the second paragraph so that the leading "{"
# and trailing "}" of the M2Doc field "{m:self.greetings.name->sep(', ')}"
# each live in their own <w:r> (matching the parser's new
# TokenIteratorFieldRewriterSplit behaviour):
#
#   "{m"                  ->  "{"  +  "m"
#   "-&gt;sep(', ')}"      ->  "->sep(', ')"  +  "}"
#
# Word's object model has no direct "split this run" verb, so we locate
# the character boundary with Find and force a run split by toggling a
# character-formatting property (Bold on, then back off) on the single
# character that must become its own run. Because the on/off round trip
# nets out to the original (unformatted) appearance, the visible content
# and formatting of the document are unchanged - only the run boundaries
# move, which is exactly what the diff calls for.

$d = $word.ActiveDocument

function Split-RunAt($rangeStart) {
    $charRange = $d.Range($rangeStart, $rangeStart + 1)
    $charRange.Bold = 1
    $charRange.Bold = 0
}

# --- "{m" -> "{" | "m" --------------------------------------------------
$finder = $d.Content
$found = $finder.Find.Execute("{m", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if ($found) {
    # $finder now spans "{m"; the split must happen right before the "m",
    # i.e. one character after the match's start.
    Split-RunAt ($finder.Start + 1)
}

# --- "->sep(', ')}" -> "->sep(', ')" | "}" ------------------------------
$finder2 = $d.Content
$found2 = $finder2.Find.Execute("->sep(', ')}", $true, $false, $false, $false, `
                                 $false, $true, 1, $false, "", 0)
if ($found2) {
    # $finder2 now spans the whole matched text; the closing "}" is its
    # last character.
    Split-RunAt ($finder2.End - 1)
}
